$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (new shared strings appended at the end:
# "Sold Stocks" / "Sold Stocks Value"), bold to match the other headers in row 5.
$ws.Range("I5").Value = "Sold Stocks"
$ws.Range("J5").Value = "Sold Stocks Value"
$ws.Range("I5:J5").Font.Bold = $true

# Size the two new columns (closest match to the bestFit sizing used by Excel).
$ws.Columns.Item(9).ColumnWidth = 9.666666666666668
$ws.Columns.Item(10).ColumnWidth = 15.666666666666668

# Update the selection to match the target workbook
$ws.Range("J6").Select()
